$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Header row: D1:M1 = Foo1..Foo10 (new shared strings), styled like C1 ---
for ($i = 0; $i -lt 10; $i++) {
    $ws3.Cells.Item(1, 4 + $i).Value = "Foo" + ($i + 1)
}
$ws3.Range("C1").Copy()
$ws3.Range("D1:M1").PasteSpecial(-4122)

# --- Data rows: D2:M32 hold (row-1); D33:M33 hold 32; D34:M34 hold 33 ---
for ($r = 2; $r -le 34; $r++) {
    $val = $r - 1
    for ($c = 4; $c -le 13; $c++) {
        $ws3.Cells.Item($r, $c).Value = $val
    }
}

# --- View: Sheet3 becomes the active/selected tab, scrolled & selected ---
$ws3.Activate() | Out-Null
$ws3.Range("D2:M34").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
